# Fix calculation of correct Excel column width
# Fix comparison of column widths if columns are part of a column group
#
# For each worksheet ("ScenarioA" and "ScenarioB") this:
#  - nudges the stored width of column A (1) and columns E:F (5:6) to the
#    corrected value (a rounding fix to the underlying width calculation)
#  - gives row 1 an explicit row height of 24pt (custom height)

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Column A (index 1): corrected width ~5.57 characters
    $ws.Columns.Item(1).ColumnWidth = 4.6667

    # Columns E and F (index 5 and 6): corrected width ~24.57 characters
    $ws.Columns.Item(5).ColumnWidth = 23.6667
    $ws.Columns.Item(6).ColumnWidth = 23.6667

    # Row 1 gains an explicit custom height of 24pt
    $ws.Rows.Item(1).RowHeight = 24
}
